$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----------------------------------------------------------------------
# Update column C (the "CP" list) for the existing rows: each list had
# contained its own Publication Number (column B) as one of the entries;
# that self-reference is removed, leaving only the *other* related codes.
# ----------------------------------------------------------------------
$ws.Cells.Item(2, 3).Value = '[''KR20090053183'']'
$ws.Cells.Item(3, 3).Value = '[''KR20210151713'']'
$ws.Cells.Item(4, 3).Value = '[''KR101923723'']'
$ws.Cells.Item(5, 3).Value = '[''JP2010097465'', ''KR20180054377'', ''US20150128062'']'
$ws.Cells.Item(6, 3).Value = '[''KR20130068593'']'
$ws.Cells.Item(7, 3).Value = '[''KR101923723'']'
$ws.Cells.Item(8, 3).Value = '[''KR102376390'']'
$ws.Cells.Item(9, 3).Value = '[''KR20140036555'']'
$ws.Cells.Item(10, 3).Value = '[''US20080147424'']'
$ws.Cells.Item(11, 3).Value = '[''KR20140036555'']'
$ws.Cells.Item(12, 3).Value = '[''KR20210063284'']'
$ws.Cells.Item(13, 3).Value = '[''US20090228550'']'
$ws.Cells.Item(14, 3).Value = '[''US20090228550'']'
$ws.Cells.Item(15, 3).Value = '[''KR102411058'']'
$ws.Cells.Item(16, 3).Value = '[''KR20120003588'']'
$ws.Cells.Item(17, 3).Value = '[''KR102341866'']'
$ws.Cells.Item(18, 3).Value = '[''US20080147424'']'
$ws.Cells.Item(19, 3).Value = '[''KR20120003588'', ''KR20210063284'']'
$ws.Cells.Item(20, 3).Value = '[''KR20140036555'']'
$ws.Cells.Item(21, 3).Value = '[''KR20140036555'', ''KR102343582'']'
$ws.Cells.Item(22, 3).Value = '[''US20080147424'']'
$ws.Cells.Item(23, 3).Value = '[''KR102343582'']'
$ws.Cells.Item(24, 3).Value = '[''KR20120003588'']'
$ws.Cells.Item(25, 3).Value = '[]'
$ws.Cells.Item(26, 3).Value = '[''KR102412142'']'
$ws.Cells.Item(27, 3).Value = '[''KR20130068593'']'
$ws.Cells.Item(28, 3).Value = '[''KR20140036555'', ''KR102343582'']'
$ws.Cells.Item(29, 3).Value = '[''KR20140036555'']'
$ws.Cells.Item(30, 3).Value = '[''KR20140036555'']'
$ws.Cells.Item(31, 3).Value = '[''KR101923723'']'
$ws.Cells.Item(32, 3).Value = '[''KR101923723'']'
$ws.Cells.Item(33, 3).Value = '[''KR20120003588'']'
$ws.Cells.Item(34, 3).Value = '[''KR20090053183'']'
$ws.Cells.Item(35, 3).Value = '[''KR20130131179'']'
$ws.Cells.Item(36, 3).Value = '[''KR20110007419'']'
$ws.Cells.Item(37, 3).Value = '[''KR20120003588'']'
$ws.Cells.Item(38, 3).Value = '[''KR20210151713'']'
$ws.Cells.Item(39, 3).Value = '[''KR100905407'']'
$ws.Cells.Item(40, 3).Value = '[''KR102343582'', ''KR101923723'']'
$ws.Cells.Item(41, 3).Value = '[''KR20090053183'']'
$ws.Cells.Item(42, 3).Value = '[''KR20210063284'']'
$ws.Cells.Item(43, 3).Value = '[''KR20120003588'', ''KR20140036555'']'
$ws.Cells.Item(44, 3).Value = '[''KR20120060597'', ''KR102404585'']'
$ws.Cells.Item(45, 3).Value = '[''KR101923723'']'
$ws.Cells.Item(46, 3).Value = '[''KR102388442'', ''KR102407595'']'
$ws.Cells.Item(47, 3).Value = '[''KR20110007419'', ''US20110014985'']'
$ws.Cells.Item(48, 3).Value = '[''US20110231434'']'
$ws.Cells.Item(49, 3).Value = '[''JP2009217387'']'
$ws.Cells.Item(50, 3).Value = '[''KR20110007419'']'
$ws.Cells.Item(51, 3).Value = '[''KR20140036555'']'
$ws.Cells.Item(52, 3).Value = '[''KR20140036555'']'
$ws.Cells.Item(53, 3).Value = '[''KR20140036555'']'
$ws.Cells.Item(54, 3).Value = '[''KR20210063284'']'
$ws.Cells.Item(55, 3).Value = '[''KR20210063284'']'
$ws.Cells.Item(56, 3).Value = '[''US8230045'', ''US8113959'', ''US20180104595'', ''US9808722'', ''US20070218987'', ''JP5159375'', ''US20090228550'']'
$ws.Cells.Item(57, 3).Value = '[''KR20140036555'', ''KR102343582'']'
$ws.Cells.Item(58, 3).Value = '[''KR20130068593'']'

# ----------------------------------------------------------------------
# The codes removed from column C above (that did not already have their
# own row) are appended as brand-new rows so every code gets one row.
# ----------------------------------------------------------------------
$ws.Cells.Item(59, 1).Value = 57
$ws.Range("A2").Copy()
$ws.Range("A59").PasteSpecial(-4122)
$ws.Cells.Item(59, 2).Value = 'KR20090053183'
$ws.Cells.Item(59, 4).Value = "'1"

$ws.Cells.Item(60, 1).Value = 58
$ws.Range("A2").Copy()
$ws.Range("A60").PasteSpecial(-4122)
$ws.Cells.Item(60, 2).Value = 'KR20210151713'
$ws.Cells.Item(60, 4).Value = "'2"

$ws.Cells.Item(61, 1).Value = 59
$ws.Range("A2").Copy()
$ws.Range("A61").PasteSpecial(-4122)
$ws.Cells.Item(61, 2).Value = 'KR101923723'
$ws.Cells.Item(61, 4).Value = "'1"

$ws.Cells.Item(62, 1).Value = 60
$ws.Range("A2").Copy()
$ws.Range("A62").PasteSpecial(-4122)
$ws.Cells.Item(62, 2).Value = 'JP2010097465'
$ws.Cells.Item(62, 4).Value = "'1"

$ws.Cells.Item(63, 1).Value = 61
$ws.Range("A2").Copy()
$ws.Range("A63").PasteSpecial(-4122)
$ws.Cells.Item(63, 2).Value = 'KR20180054377'
$ws.Cells.Item(63, 4).Value = "'1"

$ws.Cells.Item(64, 1).Value = 62
$ws.Range("A2").Copy()
$ws.Range("A64").PasteSpecial(-4122)
$ws.Cells.Item(64, 2).Value = 'US20150128062'
$ws.Cells.Item(64, 4).Value = "'1"

$ws.Cells.Item(65, 1).Value = 63
$ws.Range("A2").Copy()
$ws.Range("A65").PasteSpecial(-4122)
$ws.Cells.Item(65, 2).Value = 'KR20130068593'
$ws.Cells.Item(65, 4).Value = "'1"

$ws.Cells.Item(66, 1).Value = 64
$ws.Range("A2").Copy()
$ws.Range("A66").PasteSpecial(-4122)
$ws.Cells.Item(66, 2).Value = 'KR20140036555'
$ws.Cells.Item(66, 4).Value = "'1"

$ws.Cells.Item(67, 1).Value = 65
$ws.Range("A2").Copy()
$ws.Range("A67").PasteSpecial(-4122)
$ws.Cells.Item(67, 2).Value = 'US20080147424'
$ws.Cells.Item(67, 4).Value = "'1"

$ws.Cells.Item(68, 1).Value = 66
$ws.Range("A2").Copy()
$ws.Range("A68").PasteSpecial(-4122)
$ws.Cells.Item(68, 2).Value = 'KR20210063284'
$ws.Cells.Item(68, 4).Value = "'2"

$ws.Cells.Item(69, 1).Value = 67
$ws.Range("A2").Copy()
$ws.Range("A69").PasteSpecial(-4122)
$ws.Cells.Item(69, 2).Value = 'KR102411058'
$ws.Cells.Item(69, 4).Value = "'2"

$ws.Cells.Item(70, 1).Value = 68
$ws.Range("A2").Copy()
$ws.Range("A70").PasteSpecial(-4122)
$ws.Cells.Item(70, 2).Value = 'KR20120003588'
$ws.Cells.Item(70, 4).Value = "'1"

$ws.Cells.Item(71, 1).Value = 69
$ws.Range("A2").Copy()
$ws.Range("A71").PasteSpecial(-4122)
$ws.Cells.Item(71, 2).Value = 'KR102341866'
$ws.Cells.Item(71, 4).Value = "'2"

$ws.Cells.Item(72, 1).Value = 70
$ws.Range("A2").Copy()
$ws.Range("A72").PasteSpecial(-4122)
$ws.Cells.Item(72, 2).Value = 'KR102343582'
$ws.Cells.Item(72, 4).Value = "'2"

$ws.Cells.Item(73, 1).Value = 71
$ws.Range("A2").Copy()
$ws.Range("A73").PasteSpecial(-4122)
$ws.Cells.Item(73, 2).Value = 'KR20130131179'
$ws.Cells.Item(73, 4).Value = "'1"

$ws.Cells.Item(74, 1).Value = 72
$ws.Range("A2").Copy()
$ws.Range("A74").PasteSpecial(-4122)
$ws.Cells.Item(74, 2).Value = 'KR20110007419'
$ws.Cells.Item(74, 4).Value = "'1"

$ws.Cells.Item(75, 1).Value = 73
$ws.Range("A2").Copy()
$ws.Range("A75").PasteSpecial(-4122)
$ws.Cells.Item(75, 2).Value = 'KR100905407'
$ws.Cells.Item(75, 4).Value = "'1"

$ws.Cells.Item(76, 1).Value = 74
$ws.Range("A2").Copy()
$ws.Range("A76").PasteSpecial(-4122)
$ws.Cells.Item(76, 2).Value = 'KR20120060597'
$ws.Cells.Item(76, 4).Value = "'1"

$ws.Cells.Item(77, 1).Value = 75
$ws.Range("A2").Copy()
$ws.Range("A77").PasteSpecial(-4122)
$ws.Cells.Item(77, 2).Value = 'KR102407595'
$ws.Cells.Item(77, 4).Value = "'2"

$ws.Cells.Item(78, 1).Value = 76
$ws.Range("A2").Copy()
$ws.Range("A78").PasteSpecial(-4122)
$ws.Cells.Item(78, 2).Value = 'US20110014985'
$ws.Cells.Item(78, 4).Value = "'1"

$ws.Cells.Item(79, 1).Value = 77
$ws.Range("A2").Copy()
$ws.Range("A79").PasteSpecial(-4122)
$ws.Cells.Item(79, 2).Value = 'US20110231434'
$ws.Cells.Item(79, 4).Value = "'1"

$ws.Cells.Item(80, 1).Value = 78
$ws.Range("A2").Copy()
$ws.Range("A80").PasteSpecial(-4122)
$ws.Cells.Item(80, 2).Value = 'JP2009217387'
$ws.Cells.Item(80, 4).Value = "'1"

$ws.Cells.Item(81, 1).Value = 79
$ws.Range("A2").Copy()
$ws.Range("A81").PasteSpecial(-4122)
$ws.Cells.Item(81, 2).Value = 'US8113959'
$ws.Cells.Item(81, 4).Value = "'1"

$ws.Cells.Item(82, 1).Value = 80
$ws.Range("A2").Copy()
$ws.Range("A82").PasteSpecial(-4122)
$ws.Cells.Item(82, 2).Value = 'US20180104595'
$ws.Cells.Item(82, 4).Value = "'1"

$ws.Cells.Item(83, 1).Value = 81
$ws.Range("A2").Copy()
$ws.Range("A83").PasteSpecial(-4122)
$ws.Cells.Item(83, 2).Value = 'US20070218987'
$ws.Cells.Item(83, 4).Value = "'1"

$ws.Cells.Item(84, 1).Value = 82
$ws.Range("A2").Copy()
$ws.Range("A84").PasteSpecial(-4122)
$ws.Cells.Item(84, 2).Value = 'JP5159375'
$ws.Cells.Item(84, 4).Value = "'1"

$excel.CutCopyMode = 0

